$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")

# Row 8: replace with new event (上海·妖漫动漫展No.57&lovelive之音游ONLY（取消）)
$ws1.Range("C8").Value = "上海·妖漫动漫展No.57&lovelive之音游ONLY（取消）"
$ws1.Range("D8").Value = "漕溪北路339号 百脑汇(漕溪北路旗舰店)"
$ws1.Range("E8").Value = "2024.05.18 10:00-05.18 18:00"
$ws1.Range("F8").Value = 12
$ws1.Range("G8").Value = "不可售"
$ws1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=85325"
$ws1.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202405/XOZ6cRfV1715081814919.jpeg"

# Row 9: becomes the event previously in row 8 (上海·次元裂缝-X Anikura二次元派对)
$ws1.Range("C9").Value = "上海·次元裂缝-X Anikura二次元派对"
$ws1.Range("D9").Value = "海潮路133号B1 JUMP工坊"
$ws1.Range("E9").Value = "2024.05.18 17:00-05.18 22:00"
$ws1.Range("F9").Value = 300
$ws1.Range("G9").Value = 70
$ws1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84409"
$ws1.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202404/z38QIjBW1713260568891.jpeg"

# Remaining "想去人数" (interest count) updates on 展览
$ws1.Range("F2").Value = 1533
$ws1.Range("F4").Value = 2112
$ws1.Range("F5").Value = 7966
$ws1.Range("F6").Value = 4921
$ws1.Range("F7").Value = 7227
$ws1.Range("F13").Value = 1194
$ws1.Range("F14").Value = 202
$ws1.Range("F15").Value = 579
$ws1.Range("F16").Value = 39
$ws1.Range("F17").Value = 252
$ws1.Range("F18").Value = 5
$ws1.Range("F19").Value = 38
$ws1.Range("F20").Value = 1278
$ws1.Range("F21").Value = 1168
$ws1.Range("F24").Value = 1277
$ws1.Range("F25").Value = 60
$ws1.Range("F26").Value = 167
$ws1.Range("F28").Value = 25
$ws1.Range("F30").Value = 232
$ws1.Range("F31").Value = 1031
$ws1.Range("F34").Value = 163
$ws1.Range("F35").Value = 144
$ws1.Range("F37").Value = 563
$ws1.Range("F42").Value = 117
$ws1.Range("F43").Value = 444
$ws1.Range("F45").Value = 619
$ws1.Range("F46").Value = 171

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 34
$ws2.Range("F15").Value = 568
$ws2.Range("F17").Value = 21
$ws2.Range("F20").Value = 213
$ws2.Range("F22").Value = 147
$ws2.Range("F27").Value = 37
$ws2.Range("F28").Value = 3
$ws2.Range("F29").Value = 41
$ws2.Range("F32").Value = 889
$ws2.Range("F34").Value = 1008
$ws2.Range("F35").Value = 621
$ws2.Range("F38").Value = 128

# --- Sheet: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 867
$ws3.Range("F6").Value = 709
$ws3.Range("F7").Value = 214
$ws3.Range("F9").Value = 1805
$ws3.Range("F10").Value = 2706

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 867
$ws4.Range("F5").Value = 709
$ws4.Range("F6").Value = 7966
$ws4.Range("F7").Value = 214
$ws4.Range("F8").Value = 4921
$ws4.Range("F11").Value = 1805
$ws4.Range("F12").Value = 2706
$ws4.Range("F14").Value = 213
$ws4.Range("F15").Value = 1194
$ws4.Range("F16").Value = 147
$ws4.Range("F18").Value = 579
$ws4.Range("F19").Value = 252
$ws4.Range("F20").Value = 1278
$ws4.Range("F22").Value = 1173
$ws4.Range("F24").Value = 1277
$ws4.Range("F25").Value = 167
$ws4.Range("F27").Value = 232
$ws4.Range("F28").Value = 1031
$ws4.Range("F29").Value = 37
$ws4.Range("F30").Value = 3
$ws4.Range("F32").Value = 41
$ws4.Range("F34").Value = 163
$ws4.Range("F36").Value = 144
$ws4.Range("F39").Value = 621
$ws4.Range("F42").Value = 117
$ws4.Range("F43").Value = 619
$ws4.Range("F46").Value = 171
